$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three data rows whose values were corrected.
$ws.Range("A51").Value = "A,Double,45,High,High,16,3,3,2,436.0,2.0,nan"
$ws.Range("A67").Value = "B,Double,0,Low,High,26,3,4,3,203,10.0,nan"
$ws.Range("A77").Value = "A,Double,45,Low,High,14,4,1,1,429.0,1.0,"

# Update the view: scroll so row 58 is at the top and select A77.
$excel.ActiveWindow.ScrollRow = 58
$ws.Range("A77").Select()
